{"js": "// Apply the five text replacements described by the diff using the\n// Word JavaScript API (Office.js). Each replacement is performed via a\n// body.search() lookup (exact phrase, case-sensitive) followed by\n// insertText(..., Word.InsertLocation.replace) on the matched range.\n\nconst replacements = [\n  {\n    oldText: \"Ativa\u00e7\u00e3o: 01/01/2021\",\n    newText: \"Ativa\u00e7\u00e3o: 01/01/2024\",\n  },\n  {\n    oldText:\n      \"Abordar os princ\u00edpios b\u00e1sicos da termodin\u00e2mica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e s\u00f3lido sobre estes princ\u00edpios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodin\u00e2mica \u00e9 aplicada na pr\u00e1tica de engenharia. Enfatizar a compreens\u00e3o da termodin\u00e2mica baseada na F\u00edsica e em argumentos f\u00edsicos, buscando incentivar o entendimento mais profundo da termodin\u00e2mica.\",\n    newText:\n      \"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais, contribuindo para gerar compet\u00eancias gerais e espec\u00edficas.Abordar os princ\u00edpios b\u00e1sicos da Termodin\u00e2mica dentro do contexto de m\u00e1quinas t\u00e9rmicas.Incentivar os alunos a identificar como a termodin\u00e2mica est\u00e1 relacionada com as principais atividades humanas, com \u00eanfase na gera\u00e7\u00e3o de pot\u00eancia e refrigera\u00e7\u00e3o.Relacionar esta disciplina com outras da grade do curso, como: F\u00edsica, Recursos Naturais, Tecnologias Limpas para Gera\u00e7\u00e3o de Energia, Termodin\u00e2mica de Materiais, Sele\u00e7\u00e3o de Materiais, Fen\u00f4menos de Transporte p/ EM, dentre outras. Desenvolver nos alunos a pr\u00e1tica da busca de informa\u00e7\u00f5es t\u00e9cnicas sobre as especifica\u00e7\u00f5es de m\u00e1quinas t\u00e9rmicas e seu funcionamento. Incentivar trabalhos em grupo, com apresenta\u00e7\u00e3o de resultados.\",\n  },\n  {\n    oldText:\n      \"1. Termodin\u00e2mica e Energia. 2. Import\u00e2ncia das unidades e an\u00e1lise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 5. Propriedades de um sistema: estados termodin\u00e2micos e equil\u00edbrio. 6. Efici\u00eancia na convers\u00e3o de energia. 7. Processos e ciclos t\u00e9rmicos. 8. Termodin\u00e2mica e o meio ambiente.\",\n    newText:\n      \"1. Termodin\u00e2mica e Energia. 2. Propriedades das subst\u00e2ncias puras 3. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 4. Propriedades de um sistema: estados termodin\u00e2micos e equil\u00edbrio. 5. Efici\u00eancia na convers\u00e3o de energia. 6. Processos e ciclos t\u00e9rmicos: equipamentos, materiais e sistemas integrados. 7. Termodin\u00e2mica e o meio ambiente\",\n  },\n  {\n    oldText:\n      \"Ser\u00e3o realizadas 2 avalia\u00e7\u00f5es, com quest\u00f5es abrangendo problemas pr\u00e1ticos e conceituais. A 1a. avalia\u00e7\u00e3o ter\u00e1 peso 1 e a 2a. avalia\u00e7\u00e3o ter\u00e1 peso 2. A nota ser\u00e1 a m\u00e9dia ponderada das 2 avalia\u00e7\u00f5es.\",\n    newText:\n      \"Aulas te\u00f3ricas expositivas com recursos de m\u00eddia variados. Ser\u00e3o realizadas pelo menos duas avalia\u00e7\u00f5es escritas abrangendo problemas num\u00e9ricos e conceituais. Trabalhos em grupo abordando problemas pr\u00e1ticos tamb\u00e9m poder\u00e3o ser solicitados. Ser\u00e3o envidados esfor\u00e7os para viabilizar viagens did\u00e1ticas a plantas de gera\u00e7\u00e3o de pot\u00eancia a fim possibilitar aos alunos o contato com ciclos t\u00e9rmicos reais.\",\n  },\n  {\n    oldText:\n      \"Ser\u00e3o aplicadas duas avalia\u00e7\u00f5es escritas (P1, com peso 1 e P2, com peso 2) que compor\u00e3o a nota final (NF). A nota final ser\u00e1 calculada atrav\u00e9s da express\u00e3o: NF = (P1 + P2)/3.\",\n    newText:\n      \"Somente a nota da \u00faltima avalia\u00e7\u00e3o escrita, aplicada ao final do semestre, ter\u00e1 peso 2. As demais provas escritas ou trabalho em grupo ter\u00e3o peso 1. A nota final ser\u00e1 a m\u00e9dia ponderada dentre as avalia\u00e7\u00f5es aplicadas.\",\n  },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText.substring(0, 60));\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the five text replacements described by the diff using the Word\n# COM object model (Find/Replace on the document's main Range/Content).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{\n        Old = \"Ativa\u00e7\u00e3o: 01/01/2021\"\n        New = \"Ativa\u00e7\u00e3o: 01/01/2024\"\n    },\n    @{\n        Old = \"Abordar os princ\u00edpios b\u00e1sicos da termodin\u00e2mica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e s\u00f3lido sobre estes princ\u00edpios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodin\u00e2mica \u00e9 aplicada na pr\u00e1tica de engenharia. Enfatizar a compreens\u00e3o da termodin\u00e2mica baseada na F\u00edsica e em argumentos f\u00edsicos, buscando incentivar o entendimento mais profundo da termodin\u00e2mica.\"\n        New = \"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais, contribuindo para gerar compet\u00eancias gerais e espec\u00edficas.Abordar os princ\u00edpios b\u00e1sicos da Termodin\u00e2mica dentro do contexto de m\u00e1quinas t\u00e9rmicas.Incentivar os alunos a identificar como a termodin\u00e2mica est\u00e1 relacionada com as principais atividades humanas, com \u00eanfase na gera\u00e7\u00e3o de pot\u00eancia e refrigera\u00e7\u00e3o.Relacionar esta disciplina com outras da grade do curso, como: F\u00edsica, Recursos Naturais, Tecnologias Limpas para Gera\u00e7\u00e3o de Energia, Termodin\u00e2mica de Materiais, Sele\u00e7\u00e3o de Materiais, Fen\u00f4menos de Transporte p/ EM, dentre outras. Desenvolver nos alunos a pr\u00e1tica da busca de informa\u00e7\u00f5es t\u00e9cnicas sobre as especifica\u00e7\u00f5es de m\u00e1quinas t\u00e9rmicas e seu funcionamento. Incentivar trabalhos em grupo, com apresenta\u00e7\u00e3o de resultados.\"\n    },\n    @{\n        Old = \"1. Termodin\u00e2mica e Energia. 2. Import\u00e2ncia das unidades e an\u00e1lise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 5. Propriedades de um sistema: estados termodin\u00e2micos e equil\u00edbrio. 6. Efici\u00eancia na convers\u00e3o de energia. 7. Processos e ciclos t\u00e9rmicos. 8. Termodin\u00e2mica e o meio ambiente.\"\n        New = \"1. Termodin\u00e2mica e Energia. 2. Propriedades das subst\u00e2ncias puras 3. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 4. Propriedades de um sistema: estados termodin\u00e2micos e equil\u00edbrio. 5. Efici\u00eancia na convers\u00e3o de energia. 6. Processos e ciclos t\u00e9rmicos: equipamentos, materiais e sistemas integrados. 7. Termodin\u00e2mica e o meio ambiente\"\n    },\n    @{\n        Old = \"Ser\u00e3o realizadas 2 avalia\u00e7\u00f5es, com quest\u00f5es abrangendo problemas pr\u00e1ticos e conceituais. A 1a. avalia\u00e7\u00e3o ter\u00e1 peso 1 e a 2a. avalia\u00e7\u00e3o ter\u00e1 peso 2. A nota ser\u00e1 a m\u00e9dia ponderada das 2 avalia\u00e7\u00f5es.\"\n        New = \"Aulas te\u00f3ricas expositivas com recursos de m\u00eddia variados. Ser\u00e3o realizadas pelo menos duas avalia\u00e7\u00f5es escritas abrangendo problemas num\u00e9ricos e conceituais. Trabalhos em grupo abordando problemas pr\u00e1ticos tamb\u00e9m poder\u00e3o ser solicitados. Ser\u00e3o envidados esfor\u00e7os para viabilizar viagens did\u00e1ticas a plantas de gera\u00e7\u00e3o de pot\u00eancia a fim possibilitar aos alunos o contato com ciclos t\u00e9rmicos reais.\"\n    },\n    @{\n        Old = \"Ser\u00e3o aplicadas duas avalia\u00e7\u00f5es escritas (P1, com peso 1 e P2, com peso 2) que compor\u00e3o a nota final (NF). A nota final ser\u00e1 calculada atrav\u00e9s da express\u00e3o: NF = (P1 + P2)/3.\"\n        New = \"Somente a nota da \u00faltima avalia\u00e7\u00e3o escrita, aplicada ao final do semestre, ter\u00e1 peso 2. As demais provas escritas ou trabalho em grupo ter\u00e3o peso 1. A nota final ser\u00e1 a m\u00e9dia ponderada dentre as avalia\u00e7\u00f5es aplicadas.\"\n    }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: replacement not found for: $($pair.Old.Substring(0, [Math]::Min(60, $pair.Old.Length)))\"\n    }\n}\n"}
